# Update crypto price/volume figures per the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds text like "59.147.36" or "552.16". Plain cells whose
# text happens to look like a single-decimal number (e.g. "552.16") would get
# silently re-interpreted by Excel as a Number on assignment, so for column D we
# force the Text number format first and restore the default "Normal" style
# afterwards (keeps the cell style byte-identical to before the edit).
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" '59.147.36'
$ws.Range("E2").Value = '  -1.39%  '
Set-TextValue "D3" '2.322.00'
$ws.Range("E3").Value = '  -4.00%  '
$ws.Range("E4").Value = '  +0.04%  '
Set-TextValue "D5" '552.16'
$ws.Range("E5").Value = '  -0.18%  '
Set-TextValue "D6" '131.45'
$ws.Range("E6").Value = '  -4.19%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -4.23%  '
$ws.Range("E9").Value = '  -2.58%  '
Set-TextValue "D10" '5.55'
$ws.Range("E10").Value = '  -2.27%  '
$ws.Range("E11").Value = '  +0.83%  '
$ws.Range("E12").Value = '  -4.50%  '
Set-TextValue "D13" '23.81'
$ws.Range("E13").Value = '  -6.02%  '
Set-TextValue "D14" '2.740.96'
$ws.Range("E14").Value = '  -3.79%  '
Set-TextValue "D15" '59.147.89'
$ws.Range("E15").Value = '  -1.28%  '
$ws.Range("E16").Value = '  -2.55%  '
Set-TextValue "D17" '2.319.60'
$ws.Range("E17").Value = '  -3.74%  '
Set-TextValue "D18" '10.83'
$ws.Range("E18").Value = '  -4.11%  '
$ws.Range("E19").Value = '  -0.86%  '
Set-TextValue "D20" '317.19'
$ws.Range("E20").Value = '  -3.32%  '
$ws.Range("E21").Value = '  -2.07%  '
Set-TextValue "D22" '1.00'
$ws.Range("E22").Value = '  +0.01%  '
Set-TextValue "D23" '63.37'
$ws.Range("E23").Value = '  -4.01%  '
$ws.Range("E24").Value = '  -4.21%  '
Set-TextValue "D26" '8.32'
$ws.Range("E26").Value = '  -3.94%  '
$ws.Range("E27").Value = '  -4.90%  '
$ws.Range("E28").Value = '  +0.13%  '
Set-TextValue "D29" '170.20'
$ws.Range("E29").Value = '  +0.54%  '
Set-TextValue "D30" '0.0₃0737'
$ws.Range("E30").Value = '  -4.89%  '
$ws.Range("E31").Value = '  -3.01%  '
$ws.Range("E32").Value = '  +4.24%  '
Set-TextValue "D33" '0.391'
$ws.Range("E33").Value = '  -3.75%  '
$ws.Range("E34").Value = '  +0.01%  '
Set-TextValue "D35" '17.82'
$ws.Range("E35").Value = '  -4.20%  '
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("E37").Value = '  -2.88%  '
$ws.Range("E39").Value = '  -3.51%  '
Set-TextValue "D40" '38.53'
$ws.Range("E40").Value = '  -2.10%  '
Set-TextValue "D41" '304.34'
$ws.Range("E41").Value = '  -6.80%  '
Set-TextValue "D42" '143.49'
$ws.Range("E42").Value = '  +2.34%  '
Set-TextValue "D43" '3.45'
$ws.Range("E43").Value = '  -5.69%  '
Set-TextValue "D44" '0.0951'
Set-TextValue "D45" '0.0501'
$ws.Range("E45").Value = '  -3.11%  '
Set-TextValue "D46" '18.71'
$ws.Range("E46").Value = '  -4.96%  '
$ws.Range("E47").Value = '  -3.37%  '
Set-TextValue "D48" '0.0215'
$ws.Range("E48").Value = '  -4.07%  '
Set-TextValue "D49" '11.03'
$ws.Range("E49").Value = '  -0.04%  '
$ws.Range("E50").Value = '  +0.35%  '
Set-TextValue "D51" '0.935'
$ws.Range("E51").Value = '  -0.87%  '
